$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("C_BackLine.conf")

for ($r = 4; $r -le 58; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($val -ne $null -and $val -like "kind=*") {
        $cell.Value2 = $val -replace "^kind=", "kind%="
    }
}
